# Reorder slides ("Cambio del orden de las diapositivas").
#
# Original order:
#   1. Proyecto de ...
#   2. Tecnologías utilizadas
#   3. Requisitos de la aplicación (según la especificación previa)
#   4. Modelo entidad / relación de ASP.NET MVC
#   5. Casos de uso más relevantes (según la especificación previa)
#   6. Modelo entidad / relación de la tienda
#   7. Muchas gracias por tu atención.
#
# Target order:
#   1. Proyecto de ...
#   2. Casos de uso más relevantes (según la especificación previa)   <- moved up from 5
#   3. Requisitos de la aplicación (según la especificación previa)
#   4. Tecnologías utilizadas                                        <- pushed down from 2
#   5. Modelo entidad / relación de ASP.NET MVC                      <- pushed down from 4
#   6. Modelo entidad / relación de la tienda
#   7. Muchas gracias por tu atención.

$p = $ppt.ActivePresentation

# Move "Casos de uso más relevantes" (currently slide 5) up to position 2.
$casosDeUso = $p.Slides.Item(5)
$casosDeUso.MoveTo(2)

# "Tecnologías utilizadas" is now at position 3 (pushed down by the previous move).
# Move it to position 4, right after "Requisitos de la aplicación".
$tecnologias = $p.Slides.Item(3)
$tecnologias.MoveTo(4)
